$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as literal text
# (mirrors the original inlineStr cells, not Excel auto-number coercion).

# Row 2
$ws.Cells.Item(2, 4).Value = '27.146.79'
$ws.Cells.Item(2, 5).Value = '  -2.04%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.563.03'
$ws.Cells.Item(3, 5).Value = '  -1.68%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.02%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '206.44'
$ws.Cells.Item(5, 5).Value = '  -0.55%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.494'
$ws.Cells.Item(6, 5).Value = '  -1.86%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +0.06%  '

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '22.13'
$ws.Cells.Item(8, 5).Value = '  -0.62%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -2.04%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -0.07%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -0.72%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.786.67'
$ws.Cells.Item(12, 5).Value = '  -1.58%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.566.35'
$ws.Cells.Item(13, 5).Value = '  -1.40%  '

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '3.77'
$ws.Cells.Item(14, 5).Value = '  -2.48%  '

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.516'
$ws.Cells.Item(15, 5).Value = '  -2.82%  '

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '63.05'
$ws.Cells.Item(16, 5).Value = '  -0.68%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '27.142.89'
$ws.Cells.Item(17, 5).Value = '  -1.99%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'ShibaInu'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18, 4).Value = '0.0₃0688'
$ws.Cells.Item(18, 5).Value = '  -1.15%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'BitcoinCash'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '211.92'
$ws.Cells.Item(19, 5).Value = '  -3.74%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -1.69%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.11%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '4.11'
$ws.Cells.Item(22, 5).Value = '  -0.89%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -2.26%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  +0.56%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '152.21'
$ws.Cells.Item(25, 5).Value = '  -0.67%  '

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '6.60'
$ws.Cells.Item(26, 5).Value = '  -3.92%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '14.85'
$ws.Cells.Item(27, 5).Value = '  -2.11%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +0.07%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  -1.90%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -1.32%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '0.0464'
$ws.Cells.Item(31, 5).Value = '  -1.09%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '1.378.43'
$ws.Cells.Item(33, 5).Value = '  +0.47%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +0.40%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +0.47%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -0.23%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '0.941'
$ws.Cells.Item(37, 5).Value = '  -4.04%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -1.51%  '

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.522'
$ws.Cells.Item(39, 5).Value = '  -3.17%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.815'
$ws.Cells.Item(40, 5).Value = '  -1.34%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +0.09%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '0.989'
$ws.Cells.Item(42, 5).Value = '  +1.32%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '1.78'
$ws.Cells.Item(43, 5).Value = '  +3.20%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Aave'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '63.43'
$ws.Cells.Item(44, 5).Value = '  -1.55%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'MXToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '2.17'
$ws.Cells.Item(45, 5).Value = '  -0.14%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -0.86%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '1.696.92'
$ws.Cells.Item(47, 5).Value = '  -1.69%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '85.43'
$ws.Cells.Item(48, 5).Value = '  -2.75%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '0.0₇0997'
$ws.Cells.Item(49, 5).Value = '  -0.87%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -1.21%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +0.24%  '
